$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Rows.Item(22)
Write-Host ($r | Get-Member | Out-String)
